$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "sum-shape-area"
$ws.Range("D1").Value = "نوع کاربری"

$ws.Range("B1:G1").Select()
